$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LocationFacetMapping")

# Insert a new row at row 86, shifting existing rows 86.. down by one
$ws.Rows.Item(86).Insert()

# Populate the new row with the new location mapping
$ws.Cells.Item(86, 1).Value = "Olin Library Reserve, Circulation Desk"
$ws.Cells.Item(86, 5).Value = "Olin Library > Reserve"

# Keep the sheet view consistent with the edited region
$ws.Activate()
$ws.Range("E86").Select()
$excel.ActiveWindow.ScrollRow = 73
$excel.ActiveWindow.ScrollColumn = 1
